$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "总计" sheet: insert the new 2022-Q4 summary row above the existing
#    2022-Q3 / 2022-Q2 rows, shifting them down. Use Copy() so the existing
#    (un-styled) data cells / (styled) A-column cells keep their original
#    formatting instead of inheriting anything new.
# ---------------------------------------------------------------------------
$sheetTotal = $wb.Worksheets.Item("总计")

# Push the existing rows down one slot: row3 -> row4, row2 -> row3
$sheetTotal.Range("A3:D3").Copy($sheetTotal.Range("A4"))
$sheetTotal.Range("A2:D2").Copy($sheetTotal.Range("A3"))

# Fix up the running index in column A for the shifted rows
$sheetTotal.Range("A3").Value = 1
$sheetTotal.Range("A4").Value = 2

# Write the new 2022-Q4 summary figures into row 2
$sheetTotal.Range("A2").Value = 0
$sheetTotal.Range("B2").Value = "2022-Q4"
$sheetTotal.Range("C2").Value = 3
$sheetTotal.Range("D2").Value = 0.35

# ---------------------------------------------------------------------------
# 2. Create the new "2022-Q4" worksheet (fund holdings detail), positioned
#    right after "总计" and before "2022-Q3". Duplicate the "2022-Q3" sheet
#    so the new sheet inherits identical styles/column widths/page setup.
# ---------------------------------------------------------------------------
$sheetQ3 = $wb.Worksheets.Item("2022-Q3")
$sheetQ3.Copy($null, $sheetTotal)

$sheetQ4 = $wb.Worksheets.Item(2)
$sheetQ4.Name = "2022-Q4"

# The template sheet only had one data row; duplicate it down to rows 3 & 4
# first so every cell (notably the styled column-A index cells) gets the
# right formatting before we overwrite the values.
$sheetQ4.Range("A2:H2").Copy($sheetQ4.Range("A3"))
$sheetQ4.Range("A2:H2").Copy($sheetQ4.Range("A4"))

# Columns B, D, E, F, G hold text that looks numeric (fund code / percentages
# etc.) in this workbook - format as Text first so entering the values below
# doesn't get silently coerced into numbers (and e.g. lose the leading zero
# in the fund code).
$sheetQ4.Range("B2:B4").NumberFormat = "@"
$sheetQ4.Range("D2:G4").NumberFormat = "@"

# Row 2: 004685
$sheetQ4.Range("A2").Value = 0
$sheetQ4.Range("B2").Value = "004685"
$sheetQ4.Range("C2").Value = "金元顺安元启灵活配置混合"
$sheetQ4.Range("D2").Value = "15.29"
$sheetQ4.Range("E2").Value = "76.11"
$sheetQ4.Range("F2").Value = "1.08"
$sheetQ4.Range("G2").Value = "0.1651"
$sheetQ4.Range("H2").Value = 2

# Row 3: 001735
$sheetQ4.Range("A3").Value = 1
$sheetQ4.Range("B3").Value = "001735"
$sheetQ4.Range("C3").Value = "广发百发大数据策略成长灵活配置混合E"
$sheetQ4.Range("D3").Value = "8.88"
$sheetQ4.Range("E3").Value = "90.16"
$sheetQ4.Range("F3").Value = "1.42"
$sheetQ4.Range("G3").Value = "0.1261"
$sheetQ4.Range("H3").Value = 5

# Row 4: 001734
$sheetQ4.Range("A4").Value = 2
$sheetQ4.Range("B4").Value = "001734"
$sheetQ4.Range("C4").Value = "广发百发大数据策略成长灵活配置混合A"
$sheetQ4.Range("D4").Value = "4.29"
$sheetQ4.Range("E4").Value = "90.16"
$sheetQ4.Range("F4").Value = "1.42"
$sheetQ4.Range("G4").Value = "0.0609"
$sheetQ4.Range("H4").Value = 5

# Restore the original active tab (2022-Q2, the last sheet), since
# duplicating the Q3 sheet above switched the active sheet to the new copy.
$wb.Worksheets.Item("2022-Q2").Activate()
